{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 1: \"Click on the dropdown arrow ... downloading.\" loses its\n// trailing period, and a brand-new sub-bullet paragraph is inserted right\n// after it: \"At this point in the course, it will most likely be the\n// Food Habits campaign\".\n// ---------------------------------------------------------------------\nlet dropdownPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"Click on the dropdown arrow for the campaign you are interested in downloading.\") {\n    dropdownPara = p;\n    break;\n  }\n}\n\nif (dropdownPara) {\n  const periodRange = dropdownPara.search(\"downloading.\", { matchCase: true });\n  periodRange.load(\"text\");\n  await context.sync();\n\n  if (periodRange.items.length > 0) {\n    periodRange.items[0].insertText(\"downloading\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  const newPara = dropdownPara.insertParagraph(\n    \"At this point in the course, it will most likely be the Food Habits campaign\",\n    Word.InsertLocation.after\n  );\n  newPara.listItemOrNullObject.level = 1;\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Change 2: \"After you click import you might notice something appeared\n// in your console.\" -> the word \"import\" becomes \"Import\" and is\n// italicized.\n// ---------------------------------------------------------------------\nasync function italicizeImport(containsText) {\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n\n  let target = null;\n  for (const p of paras.items) {\n    if (p.text.indexOf(containsText) >= 0) {\n      target = p;\n      break;\n    }\n  }\n  if (!target) {\n    return;\n  }\n\n  const importRange = target.search(\"import\", { matchCase: true });\n  importRange.load(\"text\");\n  await context.sync();\n  if (importRange.items.length === 0) {\n    return;\n  }\n  importRange.items[0].insertText(\"Import\", Word.InsertLocation.replace);\n  await context.sync();\n\n  const newImportRange = target.search(\"Import\", { matchCase: true });\n  newImportRange.load(\"text\");\n  await context.sync();\n  if (newImportRange.items.length > 0) {\n    newImportRange.items[0].font.italic = true;\n    await context.sync();\n  }\n}\n\nawait italicizeImport(\"After you click import you might notice something appeared in your console.\");\n\n// ---------------------------------------------------------------------\n// Change 3: \"... uses to read your data when you clicked the import\n// button.\" -> the word \"import\" becomes \"Import\" and is italicized.\n// ---------------------------------------------------------------------\nawait italicizeImport(\"uses to read your data when you clicked the import button.\");\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1: \"Click on the dropdown arrow ... downloading.\" loses its\n# trailing period, and a brand-new sub-bullet paragraph is inserted right\n# after it: \"At this point in the course, it will most likely be the\n# Food Habits campaign\".\n# ---------------------------------------------------------------------\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n  if ($t -eq \"Click on the dropdown arrow for the campaign you are interested in downloading.`r\") {\n    $r = $p.Range\n    $r.MoveEnd(1, -1)\n    $r.Text = \"Click on the dropdown arrow for the campaign you are interested in downloading\"\n\n    $p.Range.InsertParagraphAfter()\n    $newP = $d.Paragraphs.Item($i + 1)\n    $newP.Range.Text = \"At this point in the course, it will most likely be the Food Habits campaign\"\n    $newP.Range.ListFormat.ListLevelNumber = 2\n    break\n  }\n}\n\n# ---------------------------------------------------------------------\n# Change 2: \"After you click import you might notice something appeared\n# in your console.\" -> the word \"import\" becomes \"Import\" and is\n# italicized.\n# ---------------------------------------------------------------------\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n  if ($t -like \"*After you click import you might notice something appeared in your console.*\") {\n    $pr = $p.Range\n    $pr.Find.ClearFormatting()\n    $pr.Find.Text = \"import\"\n    $pr.Find.MatchCase = $true\n    $found = $pr.Find.Execute()\n    if ($found) {\n      $pr.Text = \"Import\"\n      $pr.Font.Italic = $true\n    }\n    break\n  }\n}\n\n# ---------------------------------------------------------------------\n# Change 3: \"... uses to read your data when you clicked the import\n# button.\" -> the word \"import\" becomes \"Import\" and is italicized.\n# ---------------------------------------------------------------------\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n  if ($t -like \"*uses to read your data when you clicked the import button.*\") {\n    $pr = $p.Range\n    $pr.Find.ClearFormatting()\n    $pr.Find.Text = \"import\"\n    $pr.Find.MatchCase = $true\n    $found = $pr.Find.Execute()\n    if ($found) {\n      $pr.Text = \"Import\"\n      $pr.Font.Italic = $true\n    }\n    break\n  }\n}\n\nWrite-Output \"done\"\n"}
